$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H9").Value = 99.55556
$ws.Range("J9").Value = 0
$ws.Range("L9").Value = 0
$ws.Range("N9").ClearContents()

$ws.Range("H92").Value = 548.72
$ws.Range("I92").Value = 473.86365
$ws.Range("J92").Value = 1097.6666
$ws.Range("K92").Value = 473.86365
$ws.Range("L92").Value = 1097.6666
$ws.Range("M92").Value = 774.13635
$ws.Range("N92").Value = -3593.6666

$ws.Range("H94").Value = 1350.2307
$ws.Range("I94").Value = 1350.2307
$ws.Range("K94").Value = 1350.2307
$ws.Range("M94").Value = -899.2307000000001

$ws.Range("H98").Value = 2727.7273
$ws.Range("I98").Value = 2500.5
$ws.Range("K98").Value = 2500.5
$ws.Range("M98").Value = -1002.5

$ws.Range("H122").Value = 2727.7273
$ws.Range("I122").Value = 2500.5
$ws.Range("K122").Value = 7501.5
$ws.Range("M122").Value = -5051.5

$ws.Range("H131").Value = 3487.6316
$ws.Range("I131").Value = 2490.2727
$ws.Range("J131").Value = 4859
$ws.Range("K131").Value = 7470.8181
$ws.Range("L131").Value = 14577
$ws.Range("M131").Value = -2430.8181
$ws.Range("N131").Value = -24657

$ws.Range("H134").Value = 129999
$ws.Range("J134").Value = 129999
$ws.Range("L134").Value = 129999
$ws.Range("N134").Value = -140139

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 5498
$ws.Range("I45").Value = 5687.7085
$ws.Range("K45").Value = 5687.7085
$ws.Range("M45").Value = -5310.7085

$ws.Range("H61").Value = 15198.214
$ws.Range("I61").Value = 11689.708
$ws.Range("K61").Value = 11689.708
$ws.Range("M61").Value = -11477.708

$ws.Range("H63").Value = 2113.75
$ws.Range("I63").Value = 1531.6666
$ws.Range("K63").Value = 1531.6666
$ws.Range("M63").Value = -845.6666

$ws.Range("H66").Value = 2113.75
$ws.Range("I66").Value = 1531.6666
$ws.Range("K66").Value = 7658.333000000001
$ws.Range("M66").Value = -4226.333000000001

$ws.Range("H122").Value = 4099.4
$ws.Range("I122").Value = 4099.4
$ws.Range("K122").Value = 12298.2
$ws.Range("M122").Value = -9848.199999999999

$ws.Range("H132").Value = 4827.5186
$ws.Range("I132").Value = 2845.5652
$ws.Range("J132").Value = 16223.75
$ws.Range("K132").Value = 8536.695599999999
$ws.Range("L132").Value = 48671.25
$ws.Range("M132").Value = -6006.695599999999
$ws.Range("N132").Value = -53731.25

$ws.Range("H136").Value = 15198.214
$ws.Range("I136").Value = 11689.708
$ws.Range("K136").Value = 35069.124
$ws.Range("M136").Value = -32519.124

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 11037.744
$ws.Range("I20").Value = 13407.286
$ws.Range("K20").Value = 13407.286
$ws.Range("M20").Value = -13160.286

$ws.Range("H22").Value = 387.57144
$ws.Range("I22").Value = 387.57144
$ws.Range("K22").Value = 387.57144
$ws.Range("M22").Value = -214.57144

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 69.40000000000001
$ws.Range("I7").Value = 74.85714
$ws.Range("J7").Value = 56.666668
$ws.Range("K7").Value = 74.85714
$ws.Range("L7").Value = 56.666668
$ws.Range("M7").Value = 38.14286
$ws.Range("N7").Value = -282.666668

$ws.Range("H22").Value = 535.8182
$ws.Range("J22").Value = 919.8
$ws.Range("L22").Value = 919.8
$ws.Range("N22").Value = -1619.8

$ws.Range("H122").Value = 2661.7334
$ws.Range("I122").Value = 1443.2222
$ws.Range("K122").Value = 4329.6666
$ws.Range("M122").Value = -1879.6666

$ws.Range("H134").Value = 3829.6743
$ws.Range("I134").Value = 2580.9375
$ws.Range("J134").Value = 7462.364
$ws.Range("K134").Value = 7742.8125
$ws.Range("L134").Value = 22387.092
$ws.Range("M134").Value = -5207.8125
$ws.Range("N134").Value = -27457.092

$ws.Range("H135").Value = 65000
$ws.Range("J135").Value = 65000
$ws.Range("L135").Value = 65000
$ws.Range("N135").Value = -75140

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 1680.3462
$ws.Range("J5").Value = 2578.7693
$ws.Range("L5").Value = 7736.3079
$ws.Range("N5").Value = -7960.3079

$ws.Range("H14").Value = 2919.625
$ws.Range("I14").Value = 2919.625
$ws.Range("K14").Value = 8758.875
$ws.Range("M14").Value = -8585.875

$ws.Range("H86").Value = 411.75
$ws.Range("I86").Value = 411.75
$ws.Range("J86").Value = 0
$ws.Range("K86").Value = 1235.25
$ws.Range("L86").Value = 0
$ws.Range("M86").Value = -49.25
$ws.Range("N86").ClearContents()

$ws.Range("H89").Value = 411.75
$ws.Range("I89").Value = 411.75
$ws.Range("J89").Value = 0
$ws.Range("K89").Value = 3705.75
$ws.Range("L89").Value = 0
$ws.Range("M89").Value = 2222.25
$ws.Range("N89").ClearContents()

$ws.Range("H92").Value = 656.61536
$ws.Range("J92").Value = 797.5
$ws.Range("L92").Value = 2392.5
$ws.Range("N92").Value = -4888.5

$ws.Range("H98").Value = 1245.2941
$ws.Range("I98").Value = 1144.1818
$ws.Range("J98").Value = 1430.6666
$ws.Range("K98").Value = 3432.5454
$ws.Range("L98").Value = 4291.9998
$ws.Range("M98").Value = -1934.5454
$ws.Range("N98").Value = -7287.9998

$ws.Range("H135").Value = 1680.3462
$ws.Range("J135").Value = 2578.7693
$ws.Range("L135").Value = 23208.9237
$ws.Range("N135").Value = -28278.9237

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 2281.3635
$ws.Range("I102").Value = 2136.111
$ws.Range("K102").Value = 2136.111
$ws.Range("M102").Value = -514.1109999999999

$ws.Range("H112").Value = 0
$ws.Range("J112").Value = 0
$ws.Range("L112").Value = 0
$ws.Range("N112").ClearContents()

$ws.Range("H113").Value = 3693.5881
$ws.Range("J113").Value = 3148.1667
$ws.Range("L113").Value = 3148.1667
$ws.Range("N113").Value = -7488.1667

$ws.Range("H126").Value = 1650
$ws.Range("J126").Value = 0
$ws.Range("L126").Value = 0
$ws.Range("N126").ClearContents()

$ws.Range("H132").Value = 13698.111
$ws.Range("I132").Value = 13698.111
$ws.Range("K132").Value = 41094.333
$ws.Range("M132").Value = -38564.333

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 5306.091
$ws.Range("I132").Value = 6125.95
$ws.Range("J132").Value = 4044.7693
$ws.Range("K132").Value = 18377.85
$ws.Range("L132").Value = 12134.3079
$ws.Range("M132").Value = -15847.85
$ws.Range("N132").Value = -17194.3079

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 1243086.6
$ws.Range("I122").Value = 1855297.1
$ws.Range("K122").Value = 5565891.300000001
$ws.Range("M122").Value = -5563441.300000001

$ws.Range("H126").Value = 11951.417
$ws.Range("J126").Value = 35498.75
$ws.Range("L126").Value = 106496.25
$ws.Range("N126").Value = -111436.25
